# datos_clientes.xlsx - agregado el interfaz de acceso de clientes
# Updates row 2 (Haber/Fecha/Vencimiento/Estado) and appends a new client row (row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reference colors used by the existing conditional styling ---
# Red fill  (00FF0000) -> used for "Vencido"
# Green fill(004CE308) -> used for "Regular"
$redColor = 255
$greenColor = 582476

# ---------------------------------------------------------------
# Row 2: Franco Valentin ferraro -> Haber 2000 -> 4500, new dates,
# Estado Vencido -> Regular (style changes from red to green)
# ---------------------------------------------------------------

# Haber (H2) must stay plain text "4500" (not get reinterpreted as a
# number) while keeping its existing (red) fill style untouched.
$ws.Range("H2").Value = "'4500"
$ws.Range("H6").Copy()
$ws.Range("H2").PasteSpecial(-4122)  # xlPasteFormats - restores original style, drops quote-prefix flag
$ws.Application.CutCopyMode = $false

$ws.Range("I2").Value = "22/01/2024"
$ws.Range("J2").Value = "22/02/2024"

$ws.Range("K2").Value = "Regular"
$ws.Range("K2").Interior.Color = $greenColor

# ---------------------------------------------------------------
# Row 9: brand-new client record appended at the bottom
# ---------------------------------------------------------------
$ws.Range("A9").Value = "fre"
$ws.Range("B9").Value = "julieta"
$ws.Range("C9").Value = 12345671

# Telefono (D9) is a purely-numeric string; it must stay text (no
# cell style set in the source, same as the rest of column D).
$ws.Range("D9").Value = "'213524896"
$ws.Range("D2").Copy()
$ws.Range("D9").PasteSpecial(-4122)  # xlPasteFormats - clears quote-prefix, keeps default (no) style
$ws.Application.CutCopyMode = $false

$ws.Range("E9").Value = "gasdgsad"
$ws.Range("F9").Value = "13/10/05"
$ws.Range("G9").Value = "pase libre"

# Haber (H9) is negative but stored as text, with the green fill style.
$ws.Range("H9").Value = "'-1000"
$ws.Range("H3").Copy()
$ws.Range("H9").PasteSpecial(-4122)  # xlPasteFormats - copies green fill (same xf as H3), clears quote-prefix
$ws.Application.CutCopyMode = $false

$ws.Range("I9").Value = "22/01/2024"
$ws.Range("J9").Value = "22/02/2024"

$ws.Range("K9").Value = "Regular"
$ws.Range("K9").Interior.Color = $greenColor
